$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb2"
$ws.Cells.Item(2, 3).Value = "Rhbdl2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 45.71598933333333
$ws.Cells.Item(2, 8).Value = 137.147968
$ws.Cells.Item(2, 9).Value = 0.6549002937372808
$ws.Cells.Item(2, 10).Value = 0.6549002937372808
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 2.214002666666667
$ws.Cells.Item(2, 14).Value = 6.642008000000001
$ws.Cells.Item(2, 15).Value = 0.9692171064132561
$ws.Cells.Item(2, 16).Value = 0.9692171064132561
$ws.Cells.Item(2, 17).Value = 101.2153222933049
$ws.Cells.Item(2, 18).Value = 910.937900639744
$ws.Cells.Item(2, 19).Value = 0.6347405676852387
$ws.Cells.Item(2, 20).Value = 0.6347405676852387

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb2"
$ws.Cells.Item(3, 3).Value = "Rhbdl2"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 45.71598933333333
$ws.Cells.Item(3, 8).Value = 137.147968
$ws.Cells.Item(3, 9).Value = 0.6549002937372808
$ws.Cells.Item(3, 10).Value = 0.6549002937372808
$ws.Cells.Item(3, 11).Value = 1.0
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.06541966666666667
$ws.Cells.Item(3, 14).Value = 0.196259
$ws.Cells.Item(3, 15).Value = 0.02863856533860833
$ws.Cells.Item(3, 16).Value = 0.02863856533860832
$ws.Cells.Item(3, 17).Value = 2.990724783523556
$ws.Cells.Item(3, 18).Value = 26.916523051712
$ws.Cells.Item(3, 19).Value = 0.0187554048524689
$ws.Cells.Item(3, 20).Value = 0.0187554048524689

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb2"
$ws.Cells.Item(4, 3).Value = "Rhbdl2"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 45.71598933333333
$ws.Cells.Item(4, 8).Value = 137.147968
$ws.Cells.Item(4, 9).Value = 0.6549002937372808
$ws.Cells.Item(4, 10).Value = 0.6549002937372808
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.004898333333333334
$ws.Cells.Item(4, 14).Value = 0.014695
$ws.Cells.Item(4, 15).Value = 0.002144328248135624
$ws.Cells.Item(4, 16).Value = 0.002144328248135624
$ws.Cells.Item(4, 17).Value = 0.2239321544177778
$ws.Cells.Item(4, 18).Value = 2.01538938976
$ws.Cells.Item(4, 19).Value = 0.001404321199573169
$ws.Cells.Item(4, 20).Value = 0.001404321199573169

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efnb2"
$ws.Cells.Item(5, 3).Value = "Rhbdl2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 12.691493
$ws.Cells.Item(5, 8).Value = 38.074479
$ws.Cells.Item(5, 9).Value = 0.1818108415648851
$ws.Cells.Item(5, 10).Value = 0.1818108415648851
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 2.214002666666667
$ws.Cells.Item(5, 14).Value = 6.642008000000001
$ws.Cells.Item(5, 15).Value = 0.9692171064132561
$ws.Cells.Item(5, 16).Value = 0.9692171064132561
$ws.Cells.Item(5, 17).Value = 28.09899934598134
$ws.Cells.Item(5, 18).Value = 252.890994113832
$ws.Cells.Item(5, 19).Value = 0.1762141777760769
$ws.Cells.Item(5, 20).Value = 0.1762141777760768

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb2"
$ws.Cells.Item(6, 3).Value = "Rhbdl2"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 12.691493
$ws.Cells.Item(6, 8).Value = 38.074479
$ws.Cells.Item(6, 9).Value = 0.1818108415648851
$ws.Cells.Item(6, 10).Value = 0.1818108415648851
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.06541966666666667
$ws.Cells.Item(6, 14).Value = 0.196259
$ws.Cells.Item(6, 15).Value = 0.02863856533860833
$ws.Cells.Item(6, 16).Value = 0.02863856533860832
$ws.Cells.Item(6, 17).Value = 0.8302732415623333
$ws.Cells.Item(6, 18).Value = 7.472459174060999
$ws.Cells.Item(6, 19).Value = 0.005206801665423328
$ws.Cells.Item(6, 20).Value = 0.005206801665423326

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb2"
$ws.Cells.Item(7, 3).Value = "Rhbdl2"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 12.691493
$ws.Cells.Item(7, 8).Value = 38.074479
$ws.Cells.Item(7, 9).Value = 0.1818108415648851
$ws.Cells.Item(7, 10).Value = 0.1818108415648851
$ws.Cells.Item(7, 11).Value = 1.0
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.004898333333333334
$ws.Cells.Item(7, 14).Value = 0.014695
$ws.Cells.Item(7, 15).Value = 0.002144328248135624
$ws.Cells.Item(7, 16).Value = 0.002144328248135624
$ws.Cells.Item(7, 17).Value = 0.06216716321166667
$ws.Cells.Item(7, 18).Value = 0.5595044689049999
$ws.Cells.Item(7, 19).Value = 0.0003898621233848934
$ws.Cells.Item(7, 20).Value = 0.0003898621233848934

$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Efnb2"
$ws.Cells.Item(8, 3).Value = "Rhbdl2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 11.24784666666667
$ws.Cells.Item(8, 8).Value = 33.74354
$ws.Cells.Item(8, 9).Value = 0.161130015850732
$ws.Cells.Item(8, 10).Value = 0.161130015850732
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 2.214002666666667
$ws.Cells.Item(8, 14).Value = 6.642008000000001
$ws.Cells.Item(8, 15).Value = 0.9692171064132561
$ws.Cells.Item(8, 16).Value = 0.9692171064132561
$ws.Cells.Item(8, 17).Value = 24.90276251425778
$ws.Cells.Item(8, 18).Value = 224.12486262832
$ws.Cells.Item(8, 19).Value = 0.1561699677191685
$ws.Cells.Item(8, 20).Value = 0.1561699677191685

$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Efnb2"
$ws.Cells.Item(9, 3).Value = "Rhbdl2"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 11.24784666666667
$ws.Cells.Item(9, 8).Value = 33.74354
$ws.Cells.Item(9, 9).Value = 0.161130015850732
$ws.Cells.Item(9, 10).Value = 0.161130015850732
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.06541966666666667
$ws.Cells.Item(9, 14).Value = 0.196259
$ws.Cells.Item(9, 15).Value = 0.02863856533860833
$ws.Cells.Item(9, 16).Value = 0.02863856533860832
$ws.Cells.Item(9, 17).Value = 0.7358303796511111
$ws.Cells.Item(9, 18).Value = 6.622473416859999
$ws.Cells.Item(9, 19).Value = 0.004614532486952183
$ws.Cells.Item(9, 20).Value = 0.004614532486952183

$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Efnb2"
$ws.Cells.Item(10, 3).Value = "Rhbdl2"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 11.24784666666667
$ws.Cells.Item(10, 8).Value = 33.74354
$ws.Cells.Item(10, 9).Value = 0.161130015850732
$ws.Cells.Item(10, 10).Value = 0.161130015850732
$ws.Cells.Item(10, 11).Value = 1.0
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.004898333333333334
$ws.Cells.Item(10, 14).Value = 0.014695
$ws.Cells.Item(10, 15).Value = 0.002144328248135624
$ws.Cells.Item(10, 16).Value = 0.002144328248135624
$ws.Cells.Item(10, 17).Value = 0.05509570225555555
$ws.Cells.Item(10, 18).Value = 0.4958613202999999
$ws.Cells.Item(10, 19).Value = 0.0003455156446112654
$ws.Cells.Item(10, 20).Value = 0.0003455156446112654

$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Efnb2"
$ws.Cells.Item(11, 3).Value = "Rhbdl2"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 2.0
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.1507006666666667
$ws.Cells.Item(11, 8).Value = 0.452102
$ws.Cells.Item(11, 9).Value = 0.00215884884710222
$ws.Cells.Item(11, 10).Value = 0.00215884884710222
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 2.214002666666667
$ws.Cells.Item(11, 14).Value = 6.642008000000001
$ws.Cells.Item(11, 15).Value = 0.9692171064132561
$ws.Cells.Item(11, 16).Value = 0.9692171064132561
$ws.Cells.Item(11, 17).Value = 0.3336516778684445
$ws.Cells.Item(11, 18).Value = 3.002865100816
$ws.Cells.Item(11, 19).Value = 0.002092393232772008
$ws.Cells.Item(11, 20).Value = 0.002092393232772007

$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Efnb2"
$ws.Cells.Item(12, 3).Value = "Rhbdl2"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 2.0
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.1507006666666667
$ws.Cells.Item(12, 8).Value = 0.452102
$ws.Cells.Item(12, 9).Value = 0.00215884884710222
$ws.Cells.Item(12, 10).Value = 0.00215884884710222
$ws.Cells.Item(12, 11).Value = 1.0
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.06541966666666667
$ws.Cells.Item(12, 14).Value = 0.196259
$ws.Cells.Item(12, 15).Value = 0.02863856533860833
$ws.Cells.Item(12, 16).Value = 0.02863856533860832
$ws.Cells.Item(12, 17).Value = 0.00985878737977778
$ws.Cells.Item(12, 18).Value = 0.088729086418
$ws.Cells.Item(12, 19).Value = [double]"6.182633376391619e-05"
$ws.Cells.Item(12, 20).Value = [double]"6.182633376391616e-05"

$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Efnb2"
$ws.Cells.Item(13, 3).Value = "Rhbdl2"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2.0
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.1507006666666667
$ws.Cells.Item(13, 8).Value = 0.452102
$ws.Cells.Item(13, 9).Value = 0.00215884884710222
$ws.Cells.Item(13, 10).Value = 0.00215884884710222
$ws.Cells.Item(13, 11).Value = 1.0
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.004898333333333334
$ws.Cells.Item(13, 14).Value = 0.014695
$ws.Cells.Item(13, 15).Value = 0.002144328248135624
$ws.Cells.Item(13, 16).Value = 0.002144328248135624
$ws.Cells.Item(13, 17).Value = 0.0007381820988888889
$ws.Cells.Item(13, 18).Value = 0.00664363889
$ws.Cells.Item(13, 19).Value = [double]"4.629280566296314e-06"
$ws.Cells.Item(13, 20).Value = [double]"4.629280566296313e-06"
